$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (raw OOXML width 15.42578125 -> 16.42578125)
$ws.Columns.Item(1).ColumnWidth = 15.592447916666666

# Update cell values in column A
$ws.Cells.Item(1, 1).Value = 0.27204752559465817
$ws.Cells.Item(2, 1).Value = -0.0059999999791422454
$ws.Cells.Item(3, 1).Value = -0.0039999999791771046
$ws.Cells.Item(4, 1).Value = -0.0079999999639177588
$ws.Cells.Item(5, 1).Value = -0.002999999978262835
$ws.Cells.Item(6, 1).Value = -0.0019999999769648724
$ws.Cells.Item(7, 1).Value = -0.0099999999509989834
$ws.Cells.Item(8, 1).Value = -0.0099999999514430726
$ws.Cells.Item(9, 1).Value = -0.0019999999790925038
$ws.Cells.Item(10, 1).Value = -0.0019999999806632474
$ws.Cells.Item(11, 1).Value = -0.0029999999774457109
$ws.Cells.Item(12, 1).Value = 0.048144320495899784
$ws.Cells.Item(13, 1).Value = -0.0034999999747542532
$ws.Cells.Item(14, 1).Value = -0.0079999999595816718
$ws.Cells.Item(15, 1).Value = -0.00099999998229982623
$ws.Cells.Item(16, 1).Value = -0.0019999999787514433
$ws.Cells.Item(17, 1).Value = -0.0019999999784028333
$ws.Cells.Item(18, 1).Value = -0.0039999999717243995
$ws.Cells.Item(19, 1).Value = -0.016642373846429415
$ws.Cells.Item(20, 1).Value = -0.0039999999846376255
$ws.Cells.Item(21, 1).Value = -0.0039999999844742007
$ws.Cells.Item(22, 1).Value = -0.0039999999843525202
$ws.Cells.Item(23, 1).Value = -0.026193374963456861
$ws.Cells.Item(24, 1).Value = -0.069229998227982925
$ws.Cells.Item(25, 1).Value = -0.019999999917858169
$ws.Cells.Item(26, 1).Value = -0.0024999999745798362
$ws.Cells.Item(27, 1).Value = -0.0024999999730268563
$ws.Cells.Item(28, 1).Value = -0.0019999999683886216
$ws.Cells.Item(29, 1).Value = -0.0069999999474514851
$ws.Cells.Item(30, 1).Value = -0.059999999771263024
$ws.Cells.Item(31, 1).Value = -0.0069999999436056726
$ws.Cells.Item(32, 1).Value = -0.0099999999335498302
$ws.Cells.Item(33, 1).Value = -0.0039999999531623587
